$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new values in column C for rows 24 and 25, matching the new
# "Rewrite     -   Lily" / "Sina" shared-string entries introduced in the
# Map032 scene update.
$ws.Range("C24").Value = "Rewrite     -   Lily"
$ws.Range("C25").Value = "Sina"
